$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit after the
#    "Can two achiral functions ... flux?" bullet.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fill in the previously-empty bullet under "Uncategorized notes on the
#    primary paper" with the new Astumian note (two runs), and re-create the
#    "_GoBack" bookmark at the end of that bullet.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute("That may not be the case.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $anchor.Paragraphs(1).Next()
$xml18 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00605779" w:rsidRDefault="00605779" w:rsidP="003E3723"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">It seems Astumian wrote a 2016 paper on chemical </w:t></w:r><w:r><w:t>and optical pumping of synthetic motors (reference by not full text are in EN).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetPara.Range.InsertXML($xml18)

# ---------------------------------------------------------------------------
# 3) Re-flow the EndNote bibliography paragraph so the run split (and the
#    lastRenderedPageBreak marker) lands between "small" and "protein".
# ---------------------------------------------------------------------------
$anchor2 = $d.Content
$null = $anchor2.Find.Execute("Molecular dynamics simulations of unprecedented", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$biblioPara = $anchor2.Paragraphs(1)
$xml22 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003E3723" w:rsidRPr="003E3723" w:rsidRDefault="003E3723" w:rsidP="003E3723"><w:pPr><w:pStyle w:val="EndNoteBibliography"/><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r w:rsidRPr="003E3723"><w:tab/><w:t xml:space="preserve">Molecular dynamics simulations of unprecedented duration now can provide new insights into biomolecular mechanisms. Analysis of a 1-ms molecular dynamics simulation of the small </w:t></w:r><w:r w:rsidRPr="003E3723"><w:lastRenderedPageBreak/><w:t>protein bovine pancreatic trypsin inhibitor reveals that its main conformations have different thermodynamic profiles and that perturbation of a single geometric variable, such as a torsion angle or interresidue distance, can select for occupancy of one or another conformational state. These results establish the basis for a mechanism that we term entropy-enthalpy transduction (EET), in which the thermodynamic character of a local perturbation, such as enthalpic binding of a small molecule, is camouflaged by the thermodynamics of a global conformational change induced by the perturbation, such as a switch into a high-entropy conformational state. It is noted that EET could occur in many systems, making measured entropies and enthalpies of folding and binding unreliable indicators of actual thermodynamic driving forces. The same mechanism might also account for the high experimental variance of measured enthalpies and entropies relative to free energies in some calorimetric studies. Finally, EET may be the physical mechanism underlying many cases of entropy-enthalpy compensation.</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'
$biblioPara.Range.InsertXML($xml22)

Write-Output "Edits applied."
